$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-08 Sunday" "2025-06-09 Monday"

Replace-Text "755÷9=" "722÷7="
Replace-Text "927÷9=" "878÷8="
Replace-Text "357÷5=" "748÷9="
Replace-Text "680÷4=" "133÷2="
Replace-Text "430÷9=" "167÷9="

Replace-Text "252÷3=" "781÷2="
Replace-Text "903÷7=" "966÷2="
Replace-Text "800÷7=" "419÷6="
Replace-Text "446÷5=" "758÷2="
Replace-Text "919÷2=" "819÷9="

Replace-Text "604÷6=" "807÷5="
Replace-Text "683÷2=" "870÷3="
Replace-Text "949÷5=" "370÷3="
Replace-Text "439÷5=" "229÷6="
Replace-Text "199÷4=" "682÷7="

Replace-Text "647÷3=" "992÷7="
Replace-Text "317÷4=" "108÷6="
Replace-Text "149÷8=" "852÷8="
Replace-Text "502÷9=" "724÷8="
Replace-Text "272÷6=" "706÷9="

Replace-Text "381÷7=" "144÷2="
Replace-Text "720÷7=" "396÷2="
Replace-Text "673÷5=" "284÷7="
Replace-Text "337÷6=" "997÷3="
Replace-Text "165÷3=" "656÷3="
